$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 304 (pushes rows 304..325 down to 305..326,
# inheriting formatting/number formats from the surrounding rows, matching the diff's
# new <row r="326"> at the end and shifted data for rows 304-326).
$ws.Rows.Item(304).Insert()

# Populate the newly inserted row 304 with the weekly record added by this edit.
$ws.Range("A304").Value2 = 10
$ws.Range("B304").Value2 = "Vega Modelo de Temuco"
$ws.Range("C304").Value2 = "La Araucanía"
$ws.Range("D304").Value2 = 44585
$ws.Range("E304").Value2 = 9
$ws.Range("F304").Value2 = 100112008
$ws.Range("G304").Value2 = "Coliflor"
$ws.Range("H304").Value2 = "Sin especificar"
$ws.Range("I304").Value2 = "Primera"
$ws.Range("J304").Value2 = 450
$ws.Range("K304").Value2 = 1000
$ws.Range("L304").Value2 = 1000
$ws.Range("M304").Value2 = 1000
$ws.Range("N304").Value2 = '$/unidad'
$ws.Range("O304").Value2 = "Provincia de Cautín"
$ws.Range("P304").Value2 = 1000
$ws.Range("Q304").Value2 = 1
$ws.Range("R304").Value2 = "Hortaliza"
